# Regenerate the "K" column (column G) values for the save_data sheet.
# The workbook tracks, per observation row, a count previously derived from
# a "Strike#" source; this edit recalculates those counts ("K") and writes
# the new values (s_vals) back into column G, leaving every other column
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G), as recalculated.
$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    13 = 2
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 0
    24 = 2
    25 = 3
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 0
    31 = 3
    33 = 1
    34 = 2
    36 = 1
    37 = 1
    38 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
